$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.444.91"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.71%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.863.34"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "311.42"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -0.12%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4777"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3799"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +3.28%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07322"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9334"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.73"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +5.08%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07807"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.877.65"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.431"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.68%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.553"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "90.34"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  -0.26%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008807"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("E19").Value = "  -0.22%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "27.542.32"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.64"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("E22").Value = "  +1.05%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.69"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +0.60%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "155.82"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +0.76%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "115.24"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  -0.52%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.08881"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +3.54%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.7579"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.596"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.02%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.725"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.74%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02045"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.121"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5581"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +6.97%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05269"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.992"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.34%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "7.045"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "8.621"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +4.75%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1525"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4893"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "10.63"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.17%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "103.02"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.654"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.73%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "67.44"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.58%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06087"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.34%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.9158"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +3.17%  "
